# Apply the LinuxForHealth re-branding / regeneration edits to the
# StructureDefinition-employer workbook (gh-pages deploy of
# LinuxForHealth/alvearie-fhir-ig@80fa500).
#
# Sheet "Metadata": canonical URL, Version, Date and Publisher move from
# the old Alvearie/ibm.com identity to the new LinuxForHealth one.
# Sheet "Elements": the same canonical URL is repeated as the "Fixed
# Value" of Extension.url (Q5) and must track the Metadata change since
# both cells share the same underlying string. Separately, the
# ele-1/ext-1 invariant text that used to sit on the top-level
# "Extension" row's Constraint(s) column (AI2) is cleared - the
# regenerated IG publisher output now carries that text only on the
# child "Extension.extension" row (AI4), which already has it.

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employer"
$metadata.Range("B3").Value = "8.0.0"
$metadata.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$metadata.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employer"

# Clear AI2's text while keeping it a (blank) text cell in the same style
# - assigning "" directly would turn it into a numeric/blank cell and
# drop its style, so force text-type via a leading apostrophe then
# restore the original cell formatting from its still-populated
# neighbour (AH2, same row/style) without touching AI2's new value.
$elements.Range("AI2").Value = "'"
$elements.Range("AH2").Copy() | Out-Null
$elements.Range("AI2").PasteSpecial(-4122) | Out-Null
